# Update the "Training Dashboard" sheet with new progress as of 04-Nov-2025.
# For rows 3 through 16:
#   - Column H ("PERIOD TO EXPIRE"): decrement the existing numeric value by 1.
#   - Column I ("LAST UPDATE"): change the text from 03-Nov-2025 to 04-Nov-2025.
#     (Written as a text formula + paste-as-values so Excel does not
#     auto-convert the date-looking text into a serial date number, and so
#     the cell keeps its original style/format.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 16; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
